$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

# Clear the region-name labels (column B) in the first list block
$ws.Range("B16:B21").ClearContents()
# Clear the revenue-category labels (column B) in the second list block
$ws.Range("B24:B29").ClearContents()

# Insert a new blank row after the revenue-category list (pushes rows 30+ down by one)
$ws.Rows(30).Insert()

# Clear the region-name labels (column H) in the discount-price table (rows shifted by +1)
$ws.Range("H34:H39").ClearContents()
# Clear the trailing total/actual-performance label in that same table
$ws.Range("H40").ClearContents()

$ws.Range("B16").Select()
